# Generate Report for Handoff
# A new handoff package was produced for "b.md" (e2e\b.md). This updates the
# localization-status report so that:
#   - the Overview sheet reflects the new "Ready for handoff" status/timestamp
#     for b.md
#   - the per-locale sheets (zh-cn, de-de) record the newly generated handoff
#     file + datetime for b.md, flag it as no longer a content duplicate, and
#     surface a warning that the handback file on record is stale
#   - the "Error Detail" column is widened so the new warning text is legible

$wb = $excel.ActiveWorkbook

$newStatus        = "Ready for handoff"
$newStatusDate    = "2016-08-30 12:46:40"
$zhHandoffFile    = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhHandoffDate    = "2016-08-30 12:46:36"
$deHandoffFile    = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$deHandoffDate    = "2016-08-30 12:46:40"
$errorDetail      = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7efd7a14009b407103d16fea6882a28d104761c1/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/003d84e9984b0b37e562454eb60fa2b06825852b/e2e/b.md."

# ---------------------------------------------------------------------------
# Overview sheet: update the status + timestamp columns for the b.md row
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Range("G3").Value = $newStatusDate

# ---------------------------------------------------------------------------
# zh-cn sheet: b.md is row 3 (A3 = "b.md")
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus
# "Content Duplicate" flips to False. Assigning the literal text "False"
# would be auto-typed as a Boolean by the Value setter (losing the existing
# text representation used throughout this column), so copy the text cell
# that already holds "False" (F2) instead of typing it.
$zhcn.Range("F2").Copy()
$zhcn.Range("F3").PasteSpecial()
$zhcn.Range("G3").Value = $zhHandoffFile
$zhcn.Range("H3").Value = $zhHandoffDate
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.09

# ---------------------------------------------------------------------------
# de-de sheet: b.md is row 3 (A3 = "b.md")
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus
$dede.Range("F2").Copy()
$dede.Range("F3").PasteSpecial()
$dede.Range("G3").Value = $deHandoffFile
$dede.Range("H3").Value = $deHandoffDate
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.09
